$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 329 - this shifts the existing rows 329-387
# down to 330-388 and copies formatting (incl. the date number format on
# column D) from the surrounding rows.
$ws.Rows(329).Insert()

# Populate the newly inserted row 329 with the new weekly price record.
# Columns A,B,C,E,F,G,Q,R repeat the same constant values used throughout
# this market/product block.
$ws.Cells.Item(329, 1).Value = 11
$ws.Cells.Item(329, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(329, 3).Value = "Bíobío"
$ws.Cells.Item(329, 4).Value = 44951
$ws.Cells.Item(329, 5).Value = 8
$ws.Cells.Item(329, 6).Value = 100114001
$ws.Cells.Item(329, 7).Value = "Papa"
$ws.Cells.Item(329, 8).Value = "Patagonia"
$ws.Cells.Item(329, 9).Value = "1a (cosecha)"
$ws.Cells.Item(329, 10).Value = 5000
$ws.Cells.Item(329, 11).Value = 11000
$ws.Cells.Item(329, 12).Value = 12000
$ws.Cells.Item(329, 13).Value = 11500
$ws.Cells.Item(329, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(329, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(329, 16).Value = 460
$ws.Cells.Item(329, 17).Value = 25
$ws.Cells.Item(329, 18).Value = "Hortaliza"
